# repull data, push all data, mean calculation
# Update column F (dSF) values for rows 3-20 and row 22.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value  = 6
$ws.Range("F4").Value  = 0
$ws.Range("F5").Value  = 1
$ws.Range("F6").Value  = -5
$ws.Range("F7").Value  = -1
$ws.Range("F8").Value  = -6
$ws.Range("F9").Value  = 6
$ws.Range("F10").Value = 3
$ws.Range("F11").Value = -3
$ws.Range("F12").Value = 4
$ws.Range("F13").Value = 1
$ws.Range("F14").Value = -5
$ws.Range("F15").Value = -2
$ws.Range("F16").Value = -2
$ws.Range("F17").Value = -4
$ws.Range("F18").Value = -3
$ws.Range("F19").Value = -1
$ws.Range("F20").Value = -3
$ws.Range("F22").Value = -1
